$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()

# --- sheet view changes (scroll position + active selection) ---
$win = $excel.ActiveWindow
$win.ScrollRow = 13
$win.ScrollColumn = 11
$ws.Range("AB25").Select()

# --- value edits (rows 20-26 block) ---
$ws.Range("X20").Value = 605
$ws.Range("X25").Value = 800
$ws.Range("X26").Formula = "=600+100+4"

# --- new row 39 entry: battery (169) ---
$ws.Range("M39").Value = "battery (169)"
$ws.Range("O39").Value = 169

# --- O41 update ---
$ws.Range("O41").Value = 38

# --- O42 formula range change ---
$ws.Range("O42").Formula = "=SUM(O39:O41)"

# --- recreate shared formula blocks exactly like Excel fill would ---
$ws.Range("AA20:AC26").Formula = "=S20+W20"
$ws.Range("AA45:AA51").Formula = "=U20+Y20"
$ws.Range("AB45:AC51").Formula = "=T20+X20"

$wb.Save()
